$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 158, shifting the existing data (rows 158-225) down
# to rows 161-228.
$ws.Rows.Item(158).Insert()
$ws.Rows.Item(159).Insert()
$ws.Rows.Item(160).Insert()

# Copy static/common columns (A,B,C,E,F,G,H,I,J,K,Q,R,T) and formatting from the
# row that now holds the old first record (row 161) into the three new rows.
$ws.Range("A161:T161").Copy()
$ws.Range("A158:T158").PasteSpecial()
$ws.Range("A161:T161").Copy()
$ws.Range("A159:T159").PasteSpecial()
$ws.Range("A161:T161").Copy()
$ws.Range("A160:T160").PasteSpecial()
$excel.CutCopyMode = 0

# New data block for date 44455 (row 158: Especial, row 159: Primera, row 160: Segunda)
$ws.Range("D158").Value = 44455
$ws.Range("L158").Value = "Especial"
$ws.Range("M158").Value = 200
$ws.Range("N158").Value = 29000
$ws.Range("O158").Value = 30000
$ws.Range("P158").Value = 29500
$ws.Range("S158").Value = 4214

$ws.Range("D159").Value = 44455
$ws.Range("L159").Value = "Primera"
$ws.Range("M159").Value = 300
$ws.Range("N159").Value = 24000
$ws.Range("O159").Value = 25000
$ws.Range("P159").Value = 24500
$ws.Range("S159").Value = 3500

$ws.Range("D160").Value = 44455
$ws.Range("L160").Value = "Segunda"
$ws.Range("M160").Value = 240
$ws.Range("N160").Value = 19000
$ws.Range("O160").Value = 20000
$ws.Range("P160").Value = 19500
$ws.Range("S160").Value = 2786

Write-Output "done"
